$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")
$ws.Rows.Item(84).Insert()
$ws.Rows.Item(83).Copy()
$ws.Rows.Item(84).PasteSpecial(-4122)
